# Refactor data preparation and file handling
# - Drop the old pandas index column ("Unnamed: 0")
# - Shift "Images"/"Gender" columns left (B->A, C->B)
# - Replace the relative image paths with the absolute Windows paths
#   used by the refactored data-loading code

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A entirely; B and C shift left to A and B.
$ws.Columns.Item(1).Delete()

# New header row
$ws.Range("A1").Value = "Images"
$ws.Range("B1").Value = "Gender"

# New absolute image paths (column A), gender labels stay in column B
$base = "C:/Users/busse/Bachelorarbeit/CICD-Pipeline-Gender-Recognition/data/img_align_celeba"

$ws.Range("A2").Value = "$base\046738.jpg"
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = "$base\112331.jpg"
$ws.Range("B3").Value = 1

$ws.Range("A4").Value = "$base\031792.jpg"
$ws.Range("B4").Value = 0

$ws.Range("A5").Value = "$base\171846.jpg"
$ws.Range("B5").Value = 1
